$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.63059999999999
$ws.Range("D3").Value = -7.550199999999992
$ws.Range("A4").Value = -21.13820000000001
$ws.Range("C4").Value = -11.34609999999999
$ws.Range("D4").Value = -6.604299999999995
$ws.Range("C5").Value = -14.90540000000001
$ws.Range("E5").Value = 12.97189999999999
$ws.Range("A6").Value = -20.26969999999999
$ws.Range("C6").Value = -11.56409999999999
$ws.Range("A7").Value = -21.46050000000001
$ws.Range("A8").Value = -20.48499999999999
$ws.Range("C8").Value = -11.8857
$ws.Range("D9").Value = -7.147699999999998
$ws.Range("D11").Value = -8.127299999999995
$ws.Range("D14").Value = -6.806199999999998
$ws.Range("A16").Value = -20.46599999999999
$ws.Range("C16").Value = -12.0799
$ws.Range("D18").Value = -8.25259999999999
$ws.Range("A20").Value = -22.80040000000003
$ws.Range("E20").Value = 13.26739999999999
$ws.Range("A21").Value = -20.49389999999999
$ws.Range("C22").Value = -11.33699999999999
$ws.Range("D25").Value = -8.261799999999996

$wb.Save()
